# Update the "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" sheets:
#   F2: 650 -> 651
#   F4: 1490 -> 1491

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 651
    $ws.Range("F4").Value = 1491
}
